# Update countries & provincias Spain
# Applies the data refresh described by the commit:
#  - Inserts "Madagascar" as a new entry between "Etiopia" and
#    "Republica de Africa Central", shifting the countries that used to
#    follow it (Republica de Africa Central, Togo, Cabo Verde, Isla de
#    Man, Mauricio) down by one row with refreshed figures.
#  - Refreshes the case counters for a handful of other countries.
#  - Bumps the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 13:35"

# --- Simple numeric refreshes (country stays on the same row) --------

# Catar (row 24)
$ws.Range("B24").Value = 37097
$ws.Range("C24").Value = 1491
$ws.Range("D24").Value = 6600
$ws.Range("E24").Value = 30481
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 16

# Suiza (row 28)
$ws.Range("E28").Value = 1066
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1892

# Noruega (row 53)
$ws.Range("E53").Value = 8001
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 234

# Barein (row 55)
$ws.Range("D55").Value = 2964
$ws.Range("E55").Value = 4556

# Australia (row 57)
$ws.Range("D57").Value = 6444
$ws.Range("E57").Value = 535

# Kazajistan (row 60)
$ws.Range("D60").Value = 3649
$ws.Range("E60").Value = 3285

# Uzbekistan (row 75)
$ws.Range("D75").Value = 2366
$ws.Range("E75").Value = 501

# Senegal (row 79)
$ws.Range("B79").Value = 2714
$ws.Range("C79").Value = 97
$ws.Range("D79").Value = 1186
$ws.Range("E79").Value = 1498

# Bosnia y Herzegovina (row 80)
$ws.Range("B80").Value = 2338
$ws.Range("C80").Value = 17
$ws.Range("D80").Value = 1557
$ws.Range("E80").Value = 645
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 136

# Maldivas (row 101)
$ws.Range("B101").Value = 1186
$ws.Range("C101").Value = 43
$ws.Range("E101").Value = 1091

# --- Madagascar inserted before "Republica de Africa Central" --------
# Row 137 becomes Madagascar (new data), and the countries that used to
# occupy rows 137-141 each shift down one row (138-142), keeping their
# own latest figures.

# Row 137: Madagascar (new)
$ws.Range("A137").Value = "Madagascar"
$ws.Range("B137").Value = 371
$ws.Range("C137").Value = 45
$ws.Range("D137").Value = 131
$ws.Range("E137").Value = 238
$ws.Range("H137").Value = 2

# Row 138: Republica de Africa Central (was row 137)
$ws.Range("A138").Value = "Republica de Africa Central"
$ws.Range("B138").Value = 366
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 18
$ws.Range("E138").Value = 348
$ws.Range("H138").Value = 0

# Row 139: Togo (was row 138)
$ws.Range("A139").Value = "Togo"
$ws.Range("B139").Value = 338
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 107
$ws.Range("E139").Value = 219
$ws.Range("H139").Value = 12

# Row 140: Cabo Verde (was row 139)
$ws.Range("A140").Value = "Cabo Verde"
$ws.Range("B140").Value = 335
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 85
$ws.Range("E140").Value = 247
$ws.Range("H140").Value = 3

# Row 141: Isla de Man (was row 140)
$ws.Range("A141").Value = "Isla de Man"
$ws.Range("B141").Value = 335
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 300
$ws.Range("E141").Value = 11
$ws.Range("H141").Value = 24

# Row 142: Mauricio (was row 141)
$ws.Range("A142").Value = "Mauricio"
$ws.Range("B142").Value = 332
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 322
$ws.Range("E142").Value = 0
$ws.Range("H142").Value = 10
